$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 0.3354567892334461
$ws.Range("C7").Value = 0.4473419570207637
$ws.Range("D7").Value = 0.3498717612155226
$ws.Range("E7").Value = 0.5914995868261639
$ws.Range("F7").Value = 0.5167282329661295
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = 0.3876927802649587
$ws.Range("C8").Value = 0.5634783286452841
$ws.Range("D8").Value = 0.6514107857474178
$ws.Range("E8").Value = 0.8071002327762133
$ws.Range("F8").Value = 0.7508283629515968
$ws.Range("G8").Value = 9
